# NEW Optimization + ATP Synthase/PMF
#
# The "lb"/"ub" columns (B, C) are derived from the "base_val" column (E):
#   B = E / 100
#   C = E * 100
# Fill these formulas down for rows 12-136 (B) and 19-136 (C). Doing this with
# Range.Formula (relative reference, filled over a contiguous block) makes
# Excel record them as shared formulas, matching how the workbook was
# actually authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("lb"): =E{row}/100, rows 12-136.
# Row 12 is a lone formula; 13:76 and 77:136 are each filled as one block so
# they come out as two shared-formula groups (matching si=0 / si=2).
$ws.Range("B12").Formula = "=E12/100"
$ws.Range("B13:B76").Formula = "=E13/100"

# Column C ("ub"): =E{row}*100, rows 19-136.
# Row 19 is a lone formula; 20:83 and 84:136 are each filled as one block so
# they come out as two shared-formula groups (matching si=1 / si=3).
$ws.Range("C19").Formula = "=E19*100"
$ws.Range("C20:C83").Formula = "=E20*100"

$ws.Range("B77:B136").Formula = "=E77/100"
$ws.Range("C84:C136").Formula = "=E84*100"

# "independent" flag (column D) flips to 0 for the newly-added buffer rows.
$ws.Range("D128:D136").Value = 0

# Restore the last on-screen selection/scroll position.
$ws.Range("E118").Select()
